# Updated cryptos list: refresh Price (col D) and Volume(1h) (col E) figures.
# Numeric-looking Price values are written with a leading apostrophe so
# Excel keeps them as literal text (preserving trailing zeros / exact
# formatting) instead of silently coercing them to floating-point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.949.89"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "1.671.85"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D5").Value = "'214.76"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.0619"
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("D10").Value = "'20.14"
$ws.Range("D11").Value = "'0.0890"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "1.907.73"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "1.676.86"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "26.964.58"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'234.44"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").Value = "'8.05"
$ws.Range("E19").Value = "  +3.66%  "
$ws.Range("E20").Value = "  -0.55%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D25").Value = "'145.74"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'0.0497"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("E31").Value = "  -0.42%  "
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").Value = "1.471.10"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("D34").Value = "'3.13"
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("D35").Value = "'1.66"
$ws.Range("E35").Value = "  +1.49%  "
$ws.Range("D36").Value = "'2.42"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").Value = "'0.577"
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("D38").Value = "'0.894"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("E40").Value = "  +7.78%  "
$ws.Range("E41").Value = "  -3.95%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +2.71%  "
$ws.Range("D44").Value = "'66.72"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "1.813.54"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("D47").Value = "'90.52"
$ws.Range("E47").Value = "  +0.57%  "
$ws.Range("D48").Value = "'1.53"
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("E49").Value = "  +1.76%  "
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").Value = "'7.66"
$ws.Range("E51").Value = "  -0.37%  "
